$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.902.81"
$ws.Range("E2").Value = "'  +1.16%  "

$ws.Range("D3").Value = "'3.111.00"
$ws.Range("E3").Value = "'  +1.53%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  +0.03%  "

$ws.Range("D5").Value = "'575.79"
$ws.Range("E5").Value = "'  -0.44%  "

$ws.Range("D6").Value = "'173.58"
$ws.Range("E6").Value = "'  +3.89%  "

$ws.Range("D8").Value = "'3.106.71"
$ws.Range("E8").Value = "'  +1.53%  "

$ws.Range("E9").Value = "'  -0.43%  "

$ws.Range("D10").Value = "'6.42"
$ws.Range("E10").Value = "'  -3.43%  "

$ws.Range("E11").Value = "'  +0.04%  "

$ws.Range("E12").Value = "'  -0.79%  "

$ws.Range("E13").Value = "'  -1.31%  "

$ws.Range("D14").Value = "'37.15"
$ws.Range("E14").Value = "'  +1.26%  "

$ws.Range("E15").Value = "'  -1.01%  "

$ws.Range("D16").Value = "'3.620.99"
$ws.Range("E16").Value = "'  +1.41%  "

$ws.Range("D17").Value = "'66.865.76"
$ws.Range("E17").Value = "'  +1.12%  "

$ws.Range("D18").Value = "'7.08"
$ws.Range("E18").Value = "'  -1.29%  "

$ws.Range("D19").Value = "'3.109.26"
$ws.Range("E19").Value = "'  +1.57%  "

$ws.Range("D20").Value = "'16.27"
$ws.Range("E20").Value = "'  +0.91%  "

$ws.Range("D21").Value = "'476.15"
$ws.Range("E21").Value = "'  +2.46%  "

$ws.Range("E22").Value = "'  -0.04%  "

$ws.Range("D23").Value = "'7.77"
$ws.Range("E23").Value = "'  +4.78%  "

$ws.Range("B24").Value = "'InternetComputer(DFINITY)"
$ws.Range("C24").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'13.35"
$ws.Range("E24").Value = "'  +4.09%  "

$ws.Range("B25").Value = "'Litecoin"
$ws.Range("C25").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'83.73"
$ws.Range("E25").Value = "'  +0.89%  "

$ws.Range("E26").Value = "'  +0.53%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "'  -0.08%  "

$ws.Range("E28").Value = "'  -1.38%  "

$ws.Range("B29").Value = "'ImmutableX"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'2.42"
$ws.Range("E29").Value = "'  +0.08%  "

$ws.Range("B30").Value = "'NEARProtocol"
$ws.Range("C30").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.92"
$ws.Range("E30").Value = "'  -2.10%  "

$ws.Range("E31").Value = "'  -0.09%  "

$ws.Range("D32").Value = "'28.66"
$ws.Range("E32").Value = "'  +1.61%  "

$ws.Range("B33").Value = "'Hedera"
$ws.Range("C33").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.114"
$ws.Range("E33").Value = "'  -1.36%  "

$ws.Range("B34").Value = "'PEPE"
$ws.Range("C34").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "'0.0$([char]0x2083)0947"
$ws.Range("E34").Value = "'  -7.94%  "

$ws.Range("E35").Value = "'  +0.17%  "

$ws.Range("D36").Value = "'5.85"
$ws.Range("E36").Value = "'  -0.19%  "

$ws.Range("D37").Value = "'0.977"
$ws.Range("E37").Value = "'  -2.09%  "

$ws.Range("D38").Value = "'47.66"
$ws.Range("E38").Value = "'  -1.59%  "

$ws.Range("E39").Value = "'  +2.75%  "

$ws.Range("D40").Value = "'49.90"
$ws.Range("E40").Value = "'  +0.01%  "

$ws.Range("D41").Value = "'0.310"
$ws.Range("E41").Value = "'  -0.51%  "

$ws.Range("E42").Value = "'  -0.19%  "

$ws.Range("E43").Value = "'  -0.96%  "

$ws.Range("D44").Value = "'2.791.48"
$ws.Range("E44").Value = "'  +1.19%  "

$ws.Range("E45").Value = "'  -1.89%  "

$ws.Range("D46").Value = "'377.75"
$ws.Range("E46").Value = "'  -0.42%  "

$ws.Range("D47").Value = "'2.56"
$ws.Range("E47").Value = "'  -11.31%  "

$ws.Range("D48").Value = "'135.97"
$ws.Range("E48").Value = "'  +1.30%  "

$ws.Range("E49").Value = "'  -0.02%  "

$ws.Range("D50").Value = "'24.69"
$ws.Range("E50").Value = "'  +1.17%  "

$ws.Range("D51").Value = "'2.20"
$ws.Range("E51").Value = "'  -1.10%  "
